$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Get-ParaByText($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

function Replace-Paragraph($needle, $innerXml) {
    $p = Get-ParaByText($needle)
    $rng = $d.Range($p.Range.Start, $p.Range.End)
    $rng.InsertXML("<w:p $wns>$innerXml</w:p>")
}

function Add-ParagraphAfterLast($innerXml) {
    $last = $d.Paragraphs.Last
    $last.Range.InsertParagraphAfter()
    $newLast = $d.Paragraphs.Last
    $newLast.Range.InsertXML("<w:p $wns>$innerXml</w:p>")
}

# --- Edit 1: merge "enviar" / "/ofrecer un puesto de trabaj" / "o" runs into one run ---
Replace-Paragraph "Poster" (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='4'/></w:numPr><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>" +
    "<w:r $wns><w:rPr><w:b/><w:bCs/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>Poster</w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>: Visita</w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>n</w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'> la página web para </w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>crear oportunidades </w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>enviar/ofrecer un puesto de trabajo</w:t></w:r>"
)

# --- Edit 2: merge trailing " web para ... interactuar en ella. " + " " into one run ---
Replace-Paragraph "Poster y los Usuarios" (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:ind w:left='0'/><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>Descripción: </w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>Es el que administra la </w:t></w:r>" +
    "<w:proofErr $wns w:type='gramStart'/>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>app</w:t></w:r>" +
    "<w:proofErr $wns w:type='gramEnd'/>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'> web para que los Poster y los Usuarios puedan interactuar en ella.  </w:t></w:r>"
)

# --- Edit 3: merge " web de forma directa y " + "realiza las solicitudes..." into one run ---
Replace-Paragraph "apliquen a las mismas" (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:ind w:left='0'/><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>Descripción: Es el que i</w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>nteractúa con la </w:t></w:r>" +
    "<w:proofErr $wns w:type='gramStart'/>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>app</w:t></w:r>" +
    "<w:proofErr $wns w:type='gramEnd'/>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'> web de forma directa y realiza las solicitudes de empleo para que los usuarios apliquen a las mismas. </w:t></w:r>"
)

# --- Edit 4: append new content (11 blank paragraphs, "Asunciones y dependencias" section,
#             "Funcionalidades" section) after the final paragraph of the body ---

for ($i = 0; $i -lt 11; $i++) {
    Add-ParagraphAfterLast "<w:pPr><w:pStyle w:val='ListParagraph'/><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>"
}

Add-ParagraphAfterLast (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:rPr><w:b/><w:bCs/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr></w:pPr>" +
    "<w:r $wns><w:rPr><w:b/><w:bCs/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space='preserve'>Asunciones y dependencias </w:t></w:r>"
)

Add-ParagraphAfterLast (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>Tanto el usuario, como el poster y el </w:t></w:r>" +
    "<w:proofErr $wns w:type='spellStart'/>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>admin</w:t></w:r>" +
    "<w:proofErr $wns w:type='spellEnd'/>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'> necesita</w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>n</w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'> conexión a internet para </w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>poder </w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>acceder a la web y utilizar los servicios </w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>que brinda la misma. </w:t></w:r>"
)

Add-ParagraphAfterLast "<w:pPr><w:pStyle w:val='ListParagraph'/><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>"
Add-ParagraphAfterLast "<w:pPr><w:pStyle w:val='ListParagraph'/><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>"

Add-ParagraphAfterLast (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:rPr><w:b/><w:bCs/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr></w:pPr>" +
    "<w:r $wns><w:rPr><w:b/><w:bCs/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr><w:t>Funcionalidades</w:t></w:r>"
)

Add-ParagraphAfterLast "<w:pPr><w:pStyle w:val='ListParagraph'/><w:rPr><w:b/><w:bCs/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr></w:pPr>"

Add-ParagraphAfterLast (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>" +
    "<w:r $wns><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>Aplicación Web: &#8729;</w:t></w:r>"
)

Add-ParagraphAfterLast (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'> La aplicación debe permitir que </w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>el </w:t></w:r>" +
    "<w:proofErr $wns w:type='spellStart'/>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>admin</w:t></w:r>" +
    "<w:proofErr $wns w:type='spellEnd'/>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'> realice ajustes en la página, así como modificar las categorías disponibles. También puede editar y eliminar cualquier puesto de trabajo publicado</w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>. &#8729; Debe permitir </w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>crear cuenta para poder crear vacantes y publicarlas</w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>. &#8729; La aplicación debe ser capaz de mostrar </w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>las vacantes realizadas por el poster y realizar las solicitudes de estas</w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>. &#8729; </w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>Un usuario no podrá publicar una vacante si no está registrado</w:t></w:r>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>. </w:t></w:r>"
)

Add-ParagraphAfterLast "<w:pPr><w:pStyle w:val='ListParagraph'/><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>"

Add-ParagraphAfterLast (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>" +
    "<w:r $wns><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>Servicio Web / API</w:t></w:r>"
)

Add-ParagraphAfterLast (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>&#8729; El API debe ser el puente directo entre la base de datos y cualquier agente externo a ella. </w:t></w:r>"
)

Add-ParagraphAfterLast "<w:pPr><w:pStyle w:val='ListParagraph'/><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>"

Add-ParagraphAfterLast (
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:rPr><w:b/><w:bCs/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>" +
    "<w:r $wns><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>&#8729; El API debe manejar las peticiones y respuestas de las aplicaciones externas, proveyéndoles la data correspondiente a su petición.</w:t></w:r>"
)

# --- Edit 5: styles.xml - add <w:semiHidden/> to DefaultParagraphFont style ---
$stylesXml = $d.Styles
foreach ($s in $stylesXml) {
    if ($s.NameLocal -eq "Default Paragraph Font") {
        $s.Hidden = $true
    }
}

Write-Output "All edits applied."
